# Revert 'cards' to commit 95cda46ab8 (Jun 25)
# Adds Airbyte metadata columns around the existing cty_* columns and
# repopulates the data rows with the pre-consolidation values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reshape columns -------------------------------------------------
# Current layout: A=cty_code B=cty_iden C=cty_labe D=updated_at
# Target layout:  A=_airbyte_ab_id B=_airbyte_emitted_at C=cty_code
#                 D=cty_iden E=cty_labe F=_airbyte_additional_properties
#                 G=source_file_path H=updated_at

# Insert two columns before the current column A so cty_code.. shift to C..
$ws.Columns("A:B").Insert()

# Layout is now: C=cty_code D=cty_iden E=cty_labe F=updated_at
# Insert two more columns before the (now) updated_at column so it moves to H
$ws.Columns("F:G").Insert()

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("F1").Value = "_airbyte_additional_properties"
$ws.Range("G1").Value = "source_file_path"

# The newly inserted A/B columns don't inherit the bold/centered header
# look the rest of row 1 carries (F1/G1 already picked it up automatically
# from the column insert) - copy the format from an existing header cell
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows -----------------------------------------------------
$sourceFile = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/COMPANIES_TYPES/2024_08_06_1722929004063_0.parquet"

$rows = @(
    @{ Row=2;  Id="cf2757e9-b255-4c73-a919-f7a872ba655e"; Iden="01" },
    @{ Row=3;  Id="6fd854a8-e39e-4249-85d8-010a3ea632f3"; Iden="02" },
    @{ Row=4;  Id="78241a6f-0922-45dc-965a-1aeea5f12f6d"; Iden="03" },
    @{ Row=5;  Id="1ed2b349-513c-4d48-aef9-0606bb9e2ce4"; Iden="04" },
    @{ Row=6;  Id="d50334b4-b692-44ee-b53d-2d03f45ee11e"; Iden="05" },
    @{ Row=7;  Id="8a0357ed-4b3c-41e8-aa27-30e5103577f0"; Iden="06" },
    @{ Row=8;  Id="06bad272-65a5-47a7-ac1b-049cbdf066cc"; Iden="08" },
    @{ Row=9;  Id="8c77df47-a401-49d0-916f-c7d20927b17c"; Iden="09" },
    @{ Row=10; Id="abb28ec1-9f08-4b8e-a217-b9d2383c29d8"; Iden="07" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # A: _airbyte_ab_id (uuid string)
    $ws.Cells.Item($row, 1).Value = $r.Id

    # B: _airbyte_emitted_at (date/time, same for every row)
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = 45510.3079196875

    # D: cty_iden now holds the zero-padded text code, forced to text so
    # the leading zero survives; ClearFormats afterwards keeps it as text
    # while dropping the helper number format back to the sheet default
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.Iden
    $ws.Cells.Item($row, 4).ClearFormats()

    # F: _airbyte_additional_properties -> blank
    $ws.Cells.Item($row, 6).Value = ""

    # G: source_file_path
    $ws.Cells.Item($row, 7).Value = $sourceFile

    # H: updated_at refreshed timestamp
    $ws.Cells.Item($row, 8).Value = 45511.29476472052
}
